# Deployment.docx: "updated deployment to include python and sqlite"
#
# Splits the run "Make sure git for command line is installed." into:
#   "Make sure git for command line" + ", python, and SQLite" + " is installed."
# and plants the (singleton) _GoBack bookmark between the 2nd and 3rd pieces
# (Word only ever has one _GoBack bookmark in a document, so re-adding it
# here moves it away from its old location near "The soundboard will load
# and work.", which reproduces the removal shown in the diff.)

$d = $word.ActiveDocument

# Locate the run of text we need to split.
$full = $d.Content
$full.Find.Execute("Make sure git for command line is installed.", $true, $false, $false, $false, $false,
                    $true, 1, $false, "", 0)
$fullStart = $full.Start

$firstPart = $d.Content
$firstPart.Find.Execute("Make sure git for command line", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
$firstEnd = $firstPart.End

# Insert the new middle segment right after "line" (before " is installed.").
$insertion = ", python, and SQLite"
$insertPoint = $d.Range($firstEnd, $firstEnd)
$insertPoint.InsertBefore($insertion)

# Force a run break between "...command line" and ", python, and SQLite..."
# without altering any character formatting: toggling Bold on then back off
# splits the run at this boundary but nets out to a no-op on the rPr.
$firstRun = $d.Range($fullStart, $firstEnd)
$firstRun.Bold = 1
$firstRun.Bold = 0

# Re-plant the _GoBack bookmark right after "...SQLite" (before the space
# that precedes "is installed."). Adding it here both splits the run again
# and relocates the bookmark from its previous position in the document.
$bmPos = $firstEnd + $insertion.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
